# Updated symbol list on Mon Dec 12 09:45:03 UTC 2022 with GitHub Actions
#
# The "Price" column (D) holds plain numeric-looking text (e.g. "281.93"),
# stored as literal text in the source sheet. Force the number format to
# Text on each target cell BEFORE writing so Excel keeps the new values as
# literal strings instead of silently auto-converting them to numbers
# (NumberFormat has to be set per-cell - applying it once to a multi-area
# union only takes effect on the union's first area).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) refresh ---
$priceUpdates = [ordered]@{
    "D2"  = "282.01"
    "D3"  = "20.99"
    "D4"  = "6.239"
    "D5"  = "0.06158"
    "D6"  = "3.573"
    "D7"  = "6.561"
    "D8"  = "1.475"
    "D9"  = "0.8174"
    "D10" = "0.01383"
    "D11" = "0.1637"
    "D12" = "0.08297"
    "D13" = "0.03537"
    "D14" = "0.03150"
    "D16" = "3.725"
    "D17" = "0.001642"
    "D18" = "0.04649"
    "D19" = "0.006454"
    "D20" = "0.006186"
    "D23" = "3.814"
    "D24" = "2.337"
    "D26" = "0.1248"
    "D40" = "0.04649"
    "D41" = "0.007119"
    "D42" = "0.004746"
    "D43" = "0.1102"
    "D44" = "0.01146"
    "D45" = "0.00006211"
    "D47" = "0.9993"
    "D48" = "0.002924"
    "D49" = "0.00001899"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# --- Rows 42/43: BKEXToken and CEJI swap places (with refreshed prices) ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E43").Value = "42BKEXTokenBKK"
